$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Phase1" to "Sheet1"
$ws.Name = "Sheet1"

# Fill in the task/header labels across row 1, columns B through W
$ws.Range("B1").Value = "GS61PR / Position and Roles of an MLRS/HIMARS Crewmember"
$ws.Range("C1").Value = "GS61DD / Operate an M993A1 Carrier Vehicle During a Fire Mission "
$ws.Range("D1").Value = "GS61EM / Operate an IETM"
$ws.Range("E1").Value = "GS61DB / Perform Gunner's Duties (M270A1-M142)"
$ws.Range("F1").Value = "GS61SG / Perform Tactical Communications Using SINCGAARS"
$ws.Range("G1").Value = "GS6100 / Conduct SINCGAARS Exam and Review"
$ws.Range("H1").Value = "GS61HR / Operate AN/PRC-150 ( C ) HARRIS Radio"
$ws.Range("I1").Value = "GS61A1 / Introduction to the M270A1 MLRS - M142 HIMARS"
$ws.Range("J1").Value = "GS61MT / Perform PMCS on M993-M142 Vehicles During Fire Mission"
$ws.Range("K1").Value = "GS61RL / Perform Drivers Duties during and MLRS-HIMARS Reload"
$ws.Range("L1").Value = "GS6120 / MLRS-HIMARS Launcher Performance Exam and Criteque"
$ws.Range("M1").Value = "GS61HM / Perform PMCS on an M985A4 HEMTT and M989A1 HEMAT"
$ws.Range("N1").Value = "GH6121 / M142 (HIMARS) Maintenance Hands on Performance Examination"
$ws.Range("O1").Value = "GH61RS / Perform M142 / M985A4 Ammo Resupply Prcedures"
$ws.Range("P1").Value = "GH6111 / Conduct Hands on Performance Exam/Critique on M1084A1P2 RSV Ammunition Resupply Vehicle"
$ws.Range("Q1").Value = "GH61DD / Operate and XM1140 Carrier Vehicle during a Fire Mission"
$ws.Range("R1").Value = "GS111 / Conduct and M985A4 HEMTT Hands on Performance Examination and Critique"
$ws.Range("S1").Value = "GS61SX / Conduct FTX"
$ws.Range("T1").Value = "GS61IP / Conduct In-Processing"
$ws.Range("U1").Value = "GS61PT / Conduct Physical Readiness Training"
$ws.Range("V1").Value = "GS61FR / Conducte FTX Recovery"
$ws.Range("W1").Value = "GS61OP / Conducte Out-Processing"

# Move the active selection from J1 to H1
$ws.Range("H1").Select()
